$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the "closing / last row" border formatting that currently lives
#    on row 18 (the last data row of the table) before we touch anything.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy() | Out-Null

# Insert two fresh rows. They are inserted at row 22 (a currently blank row,
# just above the signature block) so no formatting gets auto-copied from a
# populated row; this keeps the style table clean. The net effect - two new
# empty rows between the table (row 18) and the signature block - is the
# same as inserting right after row 18, since rows 19-22 are all blank.
$ws.Rows.Item(22).Insert() | Out-Null
$ws.Rows.Item(22).Insert() | Out-Null

# Row 20 becomes the new "last" table row -> give it the border/format that
# row 18 used to have.
$ws.Range("B20:J20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Re-format row 18 (now an interior row) and the brand new row 19 with the
#    regular interior-row formatting used by rows 16/17.
# ---------------------------------------------------------------------------
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Content updates on the existing header / summary cells.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 234413
$ws.Range("C13").Value2 = 3
$ws.Range("F13").Value2 = 5

# Row 16 / 17: new "Valor Mora" total and the swapped "Periodo Mora" labels.
$ws.Range("E16").Value2 = "2307"
$ws.Range("G16").Value2 = 828116
$ws.Range("E17").Value2 = "2306"
$ws.Range("G17").Value2 = 828116

# ---------------------------------------------------------------------------
# 4. Row 18 keeps its worker (YUNAIDIS ...) - values are untouched, only the
#    style changed above, so nothing further is required here.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5. Row 19 / 20: new worker SHIRLEY EDITH BANQUETH GARCES, two periods.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1047388942"
$ws.Range("D19").Value2 = "SHIRLEY EDITH BANQUETH GARCES"
$ws.Range("E19").Value2 = "2507"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1047388942"
$ws.Range("D20").Value2 = "SHIRLEY EDITH BANQUETH GARCES"
$ws.Range("E20").Value2 = "2506"
$ws.Range("F20").Value2 = 56940
$ws.Range("G20").Value2 = 1423500

# ---------------------------------------------------------------------------
# 6. Column D best-fit width for the (now slightly different) longest name.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).AutoFit() | Out-Null

Write-Output "edit complete"
